# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 - they all share the
#    same underlying string).
# 2) Narrow the "Status"-ish columns (Overview E:F, zh-cn C, de-de C) from
#    ~17.22 chars down to ~13.41 chars of stored width. The Excel column
#    width model only stores whole-pixel widths, so 12.5 characters is the
#    closest achievable ColumnWidth (snaps to the same pixel width as the
#    target).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Update status text ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow columns ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
